$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 2374
$ws1.Range("F7").Value = 509
$ws1.Range("F16").Value = 282
$ws1.Range("F22").Value = 216
$ws1.Range("F23").Value = 3935
$ws1.Range("F24").Value = 5301
$ws1.Range("F28").Value = 3393
$ws1.Range("F34").Value = 1228
$ws1.Range("F40").Value = 43

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 1037

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2632

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2632
$ws4.Range("F8").Value = 2374
$ws4.Range("F9").Value = 509
$ws4.Range("F12").Value = 1037
$ws4.Range("F19").Value = 282
$ws4.Range("F23").Value = 3935
$ws4.Range("F24").Value = 5301
$ws4.Range("F28").Value = 3393
$ws4.Range("F33").Value = 1228
